$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the cryptos price/volume refresh.
# Numeric-looking "Price" column values need to be forced to Text so Excel
# does not auto-convert them to numbers (losing formatting / trailing zeros),
# matching how the source data is stored as plain text.

$ws.Range('D2').Value = '43.413.15'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '2.313.20'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.17%  '
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.73'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.40%  '
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('D14').Value = '2.670.46'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').Value = '2.300.16'
$ws.Range('E16').Value = '  -4.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.810'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').Value = '43.321.48'
$ws.Range('E18').Value = '  +3.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').Value = '0.0₃0934'
$ws.Range('E20').Value = '  +3.65%  '
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.15'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('E24').Value = '  +1.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.75'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.28%  '
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +7.17%  '
$ws.Range('E35').Value = '  +1.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0742'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('E38').Value = '  +4.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('E41').Value = '  +5.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.31'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('E43').Value = '  +4.64%  '
$ws.Range('D44').Value = '1.978.07'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('E46').Value = '  +3.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.83%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.96%  '
$ws.Range('D51').Value = '2.536.09'
$ws.Range('E51').Value = '  +1.78%  '
